$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 8006
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 8006
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 8006
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -8344

$ws.Range("H15").Value = 1407.3011
$ws.Range("I15").Value = 1407.3011
$ws.Range("K15").Value = 4221.9033
$ws.Range("M15").Value = -4052.9033

$ws.Range("H21").Value = 15766.223
$ws.Range("I21").Value = 4750
$ws.Range("J21").Value = 18913.715
$ws.Range("K21").Value = 4750
$ws.Range("L21").Value = 18913.715
$ws.Range("M21").Value = -4282
$ws.Range("N21").Value = -19849.715

$ws.Range("H23").Value = 15766.223
$ws.Range("I23").Value = 4750
$ws.Range("J23").Value = 18913.715
$ws.Range("K23").Value = 4750
$ws.Range("L23").Value = 18913.715
$ws.Range("M23").Value = -4516
$ws.Range("N23").Value = -19381.715

$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -26240

$ws.Range("H132").Value = 4982.391
$ws.Range("I132").Value = 4747.143
$ws.Range("J132").Value = 7452.5
$ws.Range("K132").Value = 14241.429
$ws.Range("L132").Value = 22357.5
$ws.Range("M132").Value = -11711.429
$ws.Range("N132").Value = -27417.5

$ws.Range("H135").Value = 1301.3846
$ws.Range("I135").Value = 816.2222
$ws.Range("J135").Value = 2393
$ws.Range("K135").Value = 7345.999800000001
$ws.Range("L135").Value = 21537
$ws.Range("M135").Value = -4810.999800000001
$ws.Range("N135").Value = -26607

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 12317.5
$ws.Range("J23").Value = 6164
$ws.Range("L23").Value = 6164
$ws.Range("N23").Value = -6682

$ws.Range("H44").Value = 22624.25
$ws.Range("I44").Value = 500
$ws.Range("J44").Value = 29999
$ws.Range("K44").Value = 500
$ws.Range("L44").Value = 29999
$ws.Range("M44").Value = -12
$ws.Range("N44").Value = -30975

$ws.Range("H58").Value = 30000
$ws.Range("J58").Value = 30000
$ws.Range("L58").Value = 30000
$ws.Range("N58").Value = -30860

$ws.Range("H61").Value = 3252.9614
$ws.Range("I61").Value = 2770.9285
$ws.Range("J61").Value = 3815.3333
$ws.Range("K61").Value = 2770.9285
$ws.Range("L61").Value = 3815.3333
$ws.Range("M61").Value = -2558.9285
$ws.Range("N61").Value = -4239.3333

$ws.Range("H63").Value = 8723.556
$ws.Range("I63").Value = 8358.799999999999
$ws.Range("J63").Value = 8863.846
$ws.Range("K63").Value = 8358.799999999999
$ws.Range("L63").Value = 8863.846
$ws.Range("M63").Value = -7672.799999999999
$ws.Range("N63").Value = -10235.846

$ws.Range("H66").Value = 8723.556
$ws.Range("I66").Value = 8358.799999999999
$ws.Range("J66").Value = 8863.846
$ws.Range("K66").Value = 41794
$ws.Range("L66").Value = 44319.23
$ws.Range("M66").Value = -38362
$ws.Range("N66").Value = -51183.23

$ws.Range("H80").Value = 20364
$ws.Range("J80").Value = 20364
$ws.Range("L80").Value = 20364
$ws.Range("N80").Value = -22360

$ws.Range("H83").Value = 20364
$ws.Range("J83").Value = 20364
$ws.Range("L83").Value = 61092
$ws.Range("N83").Value = -71076

$ws.Range("H123").Value = 38591.715
$ws.Range("J123").Value = 38591.715
$ws.Range("L123").Value = 38591.715
$ws.Range("N123").Value = -48391.715

$ws.Range("H136").Value = 3252.9614
$ws.Range("I136").Value = 2770.9285
$ws.Range("J136").Value = 3815.3333
$ws.Range("K136").Value = 8312.7855
$ws.Range("L136").Value = 11445.9999
$ws.Range("M136").Value = -5762.7855
$ws.Range("N136").Value = -16545.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 50000000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H82").Value = 16551.125
$ws.Range("I82").Value = 9804.666999999999
$ws.Range("J82").Value = 20599
$ws.Range("K82").Value = 9804.666999999999
$ws.Range("L82").Value = 20599
$ws.Range("M82").Value = -9421.666999999999
$ws.Range("N82").Value = -21365

$ws.Range("H85").Value = 16551.125
$ws.Range("I85").Value = 9804.666999999999
$ws.Range("J85").Value = 20599
$ws.Range("K85").Value = 9804.666999999999
$ws.Range("L85").Value = 20599
$ws.Range("M85").Value = -8478.666999999999
$ws.Range("N85").Value = -23251

$ws.Range("H107").Value = 51662.75
$ws.Range("I107").Value = 84617.086
$ws.Range("K107").Value = 84617.086
$ws.Range("M107").Value = -82697.086

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6348.825
$ws.Range("I31").Value = 1752
$ws.Range("J31").Value = 7683.387
$ws.Range("K31").Value = 1752
$ws.Range("L31").Value = 7683.387
$ws.Range("M31").Value = -1457
$ws.Range("N31").Value = -8273.386999999999

$ws.Range("H34").Value = 6348.825
$ws.Range("I34").Value = 1752
$ws.Range("J34").Value = 7683.387
$ws.Range("K34").Value = 1752
$ws.Range("L34").Value = 7683.387
$ws.Range("M34").Value = -1550
$ws.Range("N34").Value = -8087.387

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 131.35715
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 131.35715
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 788.1428999999999
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1014.1429

$ws.Range("H17").Value = 7772
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 7772
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 23316
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -23654

$ws.Range("H20").Value = 988.2759
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 987.4074000000001
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 2962.2222
$ws.Range("M20").Value = -2773
$ws.Range("N20").Value = -3416.2222

$ws.Range("H49").Value = 6260.9
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 7101.5
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 21304.5
$ws.Range("M49").Value = -14844
$ws.Range("N49").Value = -21616.5

$ws.Range("H68").Value = 5965.222
$ws.Range("I68").Value = 23575.5
$ws.Range("J68").Value = 933.7143
$ws.Range("K68").Value = 70726.5
$ws.Range("L68").Value = 2801.1429
$ws.Range("M68").Value = -69915.5
$ws.Range("N68").Value = -4423.1429

$ws.Range("H71").Value = 5965.222
$ws.Range("I71").Value = 23575.5
$ws.Range("J71").Value = 933.7143
$ws.Range("K71").Value = 212179.5
$ws.Range("L71").Value = 8403.4287
$ws.Range("M71").Value = -208123.5
$ws.Range("N71").Value = -16515.4287

$ws.Range("H92").Value = 755.7778
$ws.Range("I92").Value = 800.3333
$ws.Range("J92").Value = 666.6667
$ws.Range("K92").Value = 2400.9999
$ws.Range("L92").Value = 2000.0001
$ws.Range("M92").Value = -1152.9999
$ws.Range("N92").Value = -4496.0001

$ws.Range("H131").Value = 3192.709
$ws.Range("I131").Value = 380
$ws.Range("J131").Value = 3473.98
$ws.Range("K131").Value = 1140
$ws.Range("L131").Value = 10421.94
$ws.Range("M131").Value = 3900
$ws.Range("N131").Value = -20501.94

$ws.Range("H137").Value = 47203.96
$ws.Range("J137").Value = 59505.555
$ws.Range("L137").Value = 178516.665
$ws.Range("N137").Value = -188716.665

$ws.Range("H140").Value = 2471.1765
$ws.Range("J140").Value = 3123.75
$ws.Range("L140").Value = 9371.25
$ws.Range("N140").Value = -19731.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 209
$ws.Range("I2").Value = 272.25
$ws.Range("J2").Value = 158.4
$ws.Range("K2").Value = 272.25
$ws.Range("L2").Value = 158.4
$ws.Range("M2").Value = -159.25
$ws.Range("N2").Value = -384.4

$ws.Range("H43").Value = 13999.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 13999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 13999.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -14301.5

$ws.Range("H46").Value = 4315.737
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 4499.9443
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 4499.9443
$ws.Range("M46").Value = -844
$ws.Range("N46").Value = -4811.9443

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864

$ws.Range("H80").Value = 1177033.6
$ws.Range("I80").Value = 4501400
$ws.Range("J80").Value = 68911.5
$ws.Range("K80").Value = 4501400
$ws.Range("L80").Value = 68911.5
$ws.Range("M80").Value = -4500402
$ws.Range("N80").Value = -70907.5

$ws.Range("H83").Value = 1177033.6
$ws.Range("I83").Value = 4501400
$ws.Range("J83").Value = 68911.5
$ws.Range("K83").Value = 22507000
$ws.Range("L83").Value = 344557.5
$ws.Range("M83").Value = -22502008
$ws.Range("N83").Value = -354541.5

$ws.Range("H99").Value = 8609.643
$ws.Range("I99").Value = 3685
$ws.Range("J99").Value = 26666.666
$ws.Range("K99").Value = 3685
$ws.Range("L99").Value = 26666.666
$ws.Range("M99").Value = -1439
$ws.Range("N99").Value = -31158.666
